$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A4").Value = "sqdwQFEWGRBET"
$ws.Range("A4").Select()
